$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E22").Value = 5
$ws.Range("F22").Value = " -5 for wrong logic"
$ws.Range("E24").Value = 10

$ws.Range("F26").Select()
$excel.ActiveWindow.TopLeftCell = $ws.Range("A15")
